$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 6.009036962753781
$ws.Cells.Item(2, 4).Value2 = 9.221098254642811
$ws.Cells.Item(2, 5).Value2 = 13.69073008171404
$ws.Cells.Item(2, 6).Value2 = 32.61958045555525
$ws.Cells.Item(2, 7).Value2 = 3.650023328736792
$ws.Cells.Item(2, 9).Value2 = 19.51697929232831
$ws.Cells.Item(2, 10).Value2 = 9.911100611200659
$ws.Cells.Item(2, 11).Value2 = 12.65788962491837
$ws.Cells.Item(2, 14).Value2 = 18.21968480044169
$ws.Cells.Item(2, 15).Value2 = 24.48441608406176
$ws.Cells.Item(3, 2).Value2 = 5.884940128791502
$ws.Cells.Item(3, 4).Value2 = 9.171578512080591
$ws.Cells.Item(3, 5).Value2 = 13.63184402722153
$ws.Cells.Item(3, 6).Value2 = 32.62516691128275
$ws.Cells.Item(3, 7).Value2 = 3.652098339158004
$ws.Cells.Item(3, 9).Value2 = 19.61469403009081
$ws.Cells.Item(3, 10).Value2 = 9.91698098525014
$ws.Cells.Item(3, 11).Value2 = 12.32479067904739
$ws.Cells.Item(3, 14).Value2 = 18.27192950793475
$ws.Cells.Item(3, 15).Value2 = 24.53484181832607
$ws.Cells.Item(4, 2).Value2 = 5.808288486707063
$ws.Cells.Item(4, 4).Value2 = 9.142646857389954
$ws.Cells.Item(4, 5).Value2 = 13.59845018284453
$ws.Cells.Item(4, 6).Value2 = 32.63685818368388
$ws.Cells.Item(4, 7).Value2 = 3.653440427983274
$ws.Cells.Item(4, 9).Value2 = 19.67775075153773
$ws.Cells.Item(4, 10).Value2 = 9.922194062725561
$ws.Cells.Item(4, 11).Value2 = 12.11733276727223
$ws.Cells.Item(4, 14).Value2 = 18.30571205467432
$ws.Cells.Item(4, 15).Value2 = 24.57142078132421
$ws.Cells.Item(5, 2).Value2 = 5.776981443528987
$ws.Cells.Item(5, 4).Value2 = 9.131236794860033
$ws.Cells.Item(5, 5).Value2 = 13.58554706696581
$ws.Cells.Item(5, 6).Value2 = 32.64369898573737
$ws.Cells.Item(5, 7).Value2 = 3.654004499802967
$ws.Cells.Item(5, 9).Value2 = 19.70421828041547
$ws.Cells.Item(5, 10).Value2 = 9.924721711757288
$ws.Cells.Item(5, 11).Value2 = 12.03217837063836
$ws.Cells.Item(5, 14).Value2 = 18.31990828314701
$ws.Cells.Item(5, 15).Value2 = 24.58773628350975
$ws.Cells.Item(6, 2).Value2 = 5.771779954449761
$ws.Cells.Item(6, 4).Value2 = 9.129365363152944
$ws.Cells.Item(6, 5).Value2 = 13.58344739290767
$ws.Cells.Item(6, 6).Value2 = 32.64496026374718
$ws.Cells.Item(6, 7).Value2 = 3.654099201508515
$ws.Cells.Item(6, 9).Value2 = 19.70865984529105
$ws.Cells.Item(6, 10).Value2 = 9.925165789871
$ws.Cells.Item(6, 11).Value2 = 12.01800524501588
$ws.Cells.Item(6, 14).Value2 = 18.32229153217283
$ws.Cells.Item(6, 15).Value2 = 24.59053047540232
$ws.Cells.Item(7, 2).Value2 = 5.80786649724713
$ws.Cells.Item(7, 4).Value2 = 9.142491427569714
$ws.Cells.Item(7, 5).Value2 = 13.59827329857977
$ws.Cells.Item(7, 6).Value2 = 32.63694203544926
$ws.Cells.Item(7, 7).Value2 = 3.653447965700934
$ws.Cells.Item(7, 9).Value2 = 19.67810457575068
$ws.Cells.Item(7, 10).Value2 = 9.922226518298007
$ws.Cells.Item(7, 11).Value2 = 12.11618666089856
$ws.Cells.Item(7, 14).Value2 = 18.30590176926902
$ws.Cells.Item(7, 15).Value2 = 24.57163511636045
$ws.Cells.Item(8, 2).Value2 = 5.96636875547035
$ws.Cells.Item(8, 4).Value2 = 9.203723695369218
$ws.Cells.Item(8, 5).Value2 = 13.66985940298169
$ws.Cells.Item(8, 6).Value2 = 32.61979192529893
$ws.Cells.Item(8, 7).Value2 = 3.650724706815312
$ws.Cells.Item(8, 9).Value2 = 19.55003759092465
$ws.Cells.Item(8, 10).Value2 = 9.912795785263352
$ws.Cells.Item(8, 11).Value2 = 12.54371836867487
$ws.Cells.Item(8, 14).Value2 = 18.23734573093154
$ws.Cells.Item(8, 15).Value2 = 24.50063513685585
$ws.Cells.Item(9, 2).Value2 = 6.271729000493191
$ws.Cells.Item(9, 4).Value2 = 9.335063424626687
$ws.Cells.Item(9, 5).Value2 = 13.83165609794445
$ws.Cells.Item(9, 6).Value2 = 32.65169538617811
$ws.Cells.Item(9, 7).Value2 = 3.64592165731769
$ws.Cells.Item(9, 9).Value2 = 19.3230786780027
$ws.Cells.Item(9, 10).Value2 = 9.906999852420835
$ws.Cells.Item(9, 11).Value2 = 13.35373312683931
$ws.Cells.Item(9, 14).Value2 = 18.11638024891487
$ws.Cells.Item(9, 15).Value2 = 24.40610811527104
$ws.Cells.Item(10, 2).Value2 = 6.490397241010498
$ws.Cells.Item(10, 4).Value2 = 9.437819738857762
$ws.Cells.Item(10, 5).Value2 = 13.96285642073481
$ws.Cells.Item(10, 6).Value2 = 32.71499346425939
$ws.Cells.Item(10, 7).Value2 = 3.642716900266602
$ws.Cells.Item(10, 9).Value2 = 19.17093771462484
$ws.Cells.Item(10, 10).Value2 = 9.910453611026282
$ws.Cells.Item(10, 11).Value2 = 13.92526548779704
$ws.Cells.Item(10, 14).Value2 = 18.03565232604988
$ws.Cells.Item(10, 15).Value2 = 24.36407333881726
$ws.Cells.Item(11, 2).Value2 = 6.588184424922061
$ws.Cells.Item(11, 4).Value2 = 9.485788516215173
$ws.Cells.Item(11, 5).Value2 = 14.02504873138033
$ws.Cells.Item(11, 6).Value2 = 32.75239716560136
$ws.Cells.Item(11, 7).Value2 = 3.64132859607045
$ws.Cells.Item(11, 9).Value2 = 19.10486762702213
$ws.Cells.Item(11, 10).Value2 = 9.913690780118166
$ws.Cells.Item(11, 11).Value2 = 14.1789772590061
$ws.Cells.Item(11, 14).Value2 = 18.00068152441376
$ws.Cells.Item(11, 15).Value2 = 24.35093008097654
$ws.Cells.Item(12, 2).Value2 = 6.624936309798313
$ws.Cells.Item(12, 4).Value2 = 9.504117036854728
$ws.Cells.Item(12, 5).Value2 = 14.0489445935739
$ws.Cells.Item(12, 6).Value2 = 32.76779263672
$ws.Cells.Item(12, 7).Value2 = 3.640812827304096
$ws.Cells.Item(12, 9).Value2 = 19.08029794735785
$ws.Cells.Item(12, 10).Value2 = 9.915155208490832
$ws.Cells.Item(12, 11).Value2 = 14.27406120002758
$ws.Cells.Item(12, 14).Value2 = 17.98768998911758
$ws.Cells.Item(12, 15).Value2 = 24.3468142330545
$ws.Cells.Item(13, 2).Value2 = 6.617034038912759
$ws.Cells.Item(13, 4).Value2 = 9.500162575891851
$ws.Cells.Item(13, 5).Value2 = 14.04378310596968
$ws.Cells.Item(13, 6).Value2 = 32.76442229287441
$ws.Cells.Item(13, 7).Value2 = 3.640923465516864
$ws.Cells.Item(13, 9).Value2 = 19.08556949925159
$ws.Cells.Item(13, 10).Value2 = 9.914829223820359
$ws.Cells.Item(13, 11).Value2 = 14.25362857341067
$ws.Cells.Item(13, 14).Value2 = 17.99047679283709
$ws.Cells.Item(13, 15).Value2 = 24.34766233235446
$ws.Cells.Item(14, 2).Value2 = 6.591213791385734
$ws.Cells.Item(14, 4).Value2 = 9.487293203403887
$ws.Cells.Item(14, 5).Value2 = 14.02700783669509
$ws.Cells.Item(14, 6).Value2 = 32.75363912111705
$ws.Cells.Item(14, 7).Value2 = 3.641285964298518
$ws.Cells.Item(14, 9).Value2 = 19.10283726394269
$ws.Cells.Item(14, 10).Value2 = 9.913806484123221
$ws.Cells.Item(14, 11).Value2 = 14.18682021321772
$ws.Cells.Item(14, 14).Value2 = 17.99960767465316
$ws.Cells.Item(14, 15).Value2 = 24.35057420053801
$ws.Cells.Item(15, 2).Value2 = 6.575360925622417
$ws.Cells.Item(15, 4).Value2 = 9.479431301895202
$ws.Cells.Item(15, 5).Value2 = 14.01677694925499
$ws.Cells.Item(15, 6).Value2 = 32.74719427929764
$ws.Cells.Item(15, 7).Value2 = 3.641509299818315
$ws.Cells.Item(15, 9).Value2 = 19.1134727689979
$ws.Cells.Item(15, 10).Value2 = 9.913211065330312
$ws.Cells.Item(15, 11).Value2 = 14.1457666030764
$ws.Cells.Item(15, 14).Value2 = 18.00523328413237
$ws.Cells.Item(15, 15).Value2 = 24.35246999709862
$ws.Cells.Item(16, 2).Value2 = 6.483969598611202
$ws.Cells.Item(16, 4).Value2 = 9.434708460881481
$ws.Cells.Item(16, 5).Value2 = 13.95884115415939
$ws.Cells.Item(16, 6).Value2 = 32.71272175908412
$ws.Cells.Item(16, 7).Value2 = 3.642809024556474
$ws.Cells.Item(16, 9).Value2 = 19.17531856909317
$ws.Cells.Item(16, 10).Value2 = 9.910275497645319
$ws.Cells.Item(16, 11).Value2 = 13.90855073326829
$ws.Cells.Item(16, 14).Value2 = 18.03797294297608
$ws.Cells.Item(16, 15).Value2 = 24.36505272473833
$ws.Cells.Item(17, 2).Value2 = 6.427446103900978
$ws.Cells.Item(17, 4).Value2 = 9.407577395803992
$ws.Cells.Item(17, 5).Value2 = 13.92393092668148
$ws.Cells.Item(17, 6).Value2 = 32.69377480737349
$ws.Cells.Item(17, 7).Value2 = 3.64362414270446
$ws.Cells.Item(17, 9).Value2 = 19.21406177417694
$ws.Cells.Item(17, 10).Value2 = 9.908900671087835
$ws.Cells.Item(17, 11).Value2 = 13.76135261805124
$ws.Cells.Item(17, 14).Value2 = 18.05850594004784
$ws.Cells.Item(17, 15).Value2 = 24.3743042890305
$ws.Cells.Item(18, 2).Value2 = 6.394778560432638
$ws.Cells.Item(18, 4).Value2 = 9.392088326341288
$ws.Cells.Item(18, 5).Value2 = 13.90408859465858
$ws.Cells.Item(18, 6).Value2 = 32.68368782475822
$ws.Cells.Item(18, 7).Value2 = 3.644099527002617
$ws.Cells.Item(18, 9).Value2 = 19.23664146387151
$ws.Cells.Item(18, 10).Value2 = 9.908266790292322
$ws.Cells.Item(18, 11).Value2 = 13.67610169417628
$ws.Cells.Item(18, 14).Value2 = 18.07048100490119
$ws.Cells.Item(18, 15).Value2 = 24.38018821318046
$ws.Cells.Item(19, 2).Value2 = 6.383692095079711
$ws.Cells.Item(19, 4).Value2 = 9.386864286767251
$ws.Cells.Item(19, 5).Value2 = 13.89741151537853
$ws.Cells.Item(19, 6).Value2 = 32.68041197679404
$ws.Cells.Item(19, 7).Value2 = 3.644261610444298
$ws.Cells.Item(19, 9).Value2 = 19.24433738879694
$ws.Cells.Item(19, 10).Value2 = 9.908079141983286
$ws.Cells.Item(19, 11).Value2 = 13.64713928091403
$ws.Cells.Item(19, 14).Value2 = 18.07456392732025
$ws.Cells.Item(19, 15).Value2 = 24.38227699462762
$ws.Cells.Item(20, 2).Value2 = 6.433479611646038
$ws.Cells.Item(20, 4).Value2 = 9.410453626409769
$ws.Cells.Item(20, 5).Value2 = 13.92762275246041
$ws.Cells.Item(20, 6).Value2 = 32.69570787019171
$ws.Cells.Item(20, 7).Value2 = 3.643536694460884
$ws.Cells.Item(20, 9).Value2 = 19.20990691140207
$ws.Cells.Item(20, 10).Value2 = 9.909030793669526
$ws.Cells.Item(20, 11).Value2 = 13.77708346748613
$ws.Cells.Item(20, 14).Value2 = 18.05630309419177
$ws.Cells.Item(20, 15).Value2 = 24.37326119772963
$ws.Cells.Item(21, 2).Value2 = 6.598805640438743
$ws.Cells.Item(21, 4).Value2 = 9.49106890413597
$ws.Cells.Item(21, 5).Value2 = 14.03192590564395
$ws.Cells.Item(21, 6).Value2 = 32.7567730358551
$ws.Cells.Item(21, 7).Value2 = 3.641179219925047
$ws.Cells.Item(21, 9).Value2 = 19.09775311466731
$ws.Cells.Item(21, 10).Value2 = 9.914100421062724
$ws.Cells.Item(21, 11).Value2 = 14.20647101117877
$ws.Cells.Item(21, 14).Value2 = 17.99691890426515
$ws.Cells.Item(21, 15).Value2 = 24.34969553123706
$ws.Cells.Item(22, 2).Value2 = 6.705220534946295
$ws.Cells.Item(22, 4).Value2 = 9.54470444856659
$ws.Cells.Item(22, 5).Value2 = 14.10209716561313
$ws.Cells.Item(22, 6).Value2 = 32.80385681490897
$ws.Cells.Item(22, 7).Value2 = 3.639696458352345
$ws.Cells.Item(22, 9).Value2 = 19.02707406396088
$ws.Cells.Item(22, 10).Value2 = 9.918803803737386
$ws.Cells.Item(22, 11).Value2 = 14.48128779274179
$ws.Cells.Item(22, 14).Value2 = 17.95957138573777
$ws.Cells.Item(22, 15).Value2 = 24.33931422958143
$ws.Cells.Item(23, 2).Value2 = 6.648585661658267
$ws.Cells.Item(23, 4).Value2 = 9.515995455864648
$ws.Cells.Item(23, 5).Value2 = 14.06446751207075
$ws.Cells.Item(23, 6).Value2 = 32.77807337209619
$ws.Cells.Item(23, 7).Value2 = 3.640482547005396
$ws.Cells.Item(23, 9).Value2 = 19.06455768240541
$ws.Cells.Item(23, 10).Value2 = 9.916166688842129
$ws.Cells.Item(23, 11).Value2 = 14.33517182483887
$ws.Cells.Item(23, 14).Value2 = 17.97937084107154
$ws.Cells.Item(23, 15).Value2 = 24.34439519247614
$ws.Cells.Item(24, 2).Value2 = 6.430752394053216
$ws.Cells.Item(24, 4).Value2 = 9.409152942292117
$ws.Cells.Item(24, 5).Value2 = 13.92595296658423
$ws.Cells.Item(24, 6).Value2 = 32.69483142174423
$ws.Cells.Item(24, 7).Value2 = 3.64357620874434
$ws.Cells.Item(24, 9).Value2 = 19.21178437246989
$ws.Cells.Item(24, 10).Value2 = 9.908971477616998
$ws.Cells.Item(24, 11).Value2 = 13.76997348967223
$ws.Cells.Item(24, 14).Value2 = 18.05729847002869
$ws.Cells.Item(24, 15).Value2 = 24.37373101921148
$ws.Cells.Item(25, 2).Value2 = 6.189955951628309
$ws.Cells.Item(25, 4).Value2 = 9.298387537334072
$ws.Cells.Item(25, 5).Value2 = 13.78566496946327
$ws.Cells.Item(25, 6).Value2 = 32.63605305708639
$ws.Cells.Item(25, 7).Value2 = 3.647163854083816
$ws.Cells.Item(25, 9).Value2 = 19.38190246897
$ws.Cells.Item(25, 10).Value2 = 9.90721090249189
$ws.Cells.Item(25, 11).Value2 = 13.13832009710742
$ws.Cells.Item(25, 14).Value2 = 18.14766935439293
$ws.Cells.Item(25, 15).Value2 = 24.42687610433381
